$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

    $ws.Range("G2").Value = 'System, system, backup@backdoor.com'
    $ws.Range("G3").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G4").Value = 'System, backup@backdoor.com'
    $ws.Range("G5").Value = 'System, backup@backdoor.com'
    $ws.Range("G6").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G10").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G11").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G12").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G13").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G14").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G15").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G17").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G29").Value = 'System, system, backup@backdoor.com'
    $ws.Range("G30").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G31").Value = 'System, backup@backdoor.com'
    $ws.Range("G32").Value = 'System, backup@backdoor.com'
    $ws.Range("G33").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G37").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G38").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G39").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G40").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G41").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G42").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G44").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G56").Value = 'System, system, backup@backdoor.com'
    $ws.Range("G57").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G58").Value = 'System, backup@backdoor.com'
    $ws.Range("G59").Value = 'System, backup@backdoor.com'
    $ws.Range("G60").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G64").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G65").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G66").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G67").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G68").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G69").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G71").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G84").Value = 'System, backup@backdoor.com'
    $ws.Range("G85").Value = 'System, backup@backdoor.com'
    $ws.Range("G86").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G87").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G88").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G89").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G90").Value = 'dnasr281@gmail.com, admin@admin.com'
    $ws.Range("G93").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G95").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G96").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G110").Value = 'System, backup@backdoor.com'
    $ws.Range("G111").Value = 'System, backup@backdoor.com'
    $ws.Range("G112").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G113").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G114").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G115").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G116").Value = 'dnasr281@gmail.com, admin@admin.com'
    $ws.Range("G119").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G121").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G122").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G136").Value = 'System, backup@backdoor.com'
    $ws.Range("G137").Value = 'System, backup@backdoor.com'
    $ws.Range("G138").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G139").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G140").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G141").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G142").Value = 'dnasr281@gmail.com, admin@admin.com'
    $ws.Range("G145").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G147").Value = 'dnasr281@gmail.com, System'
    $ws.Range("G148").Value = 'dnasr281@gmail.com, System'
